$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2): title change only, Status (B2) stays "Added"
$ws.Range("A2").Value = "Nomads of Driftland"

# Rows 3-16 (A3:A16): replace the list of game titles, Status column (B) unchanged
$newTitles = @(
    "Slapshot: Rebound",
    "Leaf Blower Revolution - Idle Game",
    "Summerland",
    "Ruff Night At The Gallery",
    "Rubber Bandits: Christmas Prologue",
    "The Life and Suffering of Sir Brante — Chapter 1&2",
    "AOD: Art Of Defense",
    "Mrs. Santa's Gift Hunt",
    "Ninja Hanrei",
    "*NEW* SCUFFED EPIC BHOP SIMULATOR 2023 (POG CHAMP)",
    "DreamWatcher",
    "From the Shadows",
    "KAKU: Ancient Seal (Alpha)",
    "Night Reverie: Prologue"
)

$row = 3
foreach ($title in $newTitles) {
    $ws.Cells.Item($row, 1).Value = $title
    $row++
}
